$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted as row 150 in the "Alcachofa" sheet,
# pushing the existing rows 150:230 down to 151:231.
$ws.Rows("150:150").Insert()

# Populate the newly inserted row 150 with the new record's data.
$ws.Range("A150").Value2 = 10
$ws.Range("B150").Value2 = "Vega Modelo de Temuco"
$ws.Range("C150").Value2 = "La Araucanía"
$ws.Range("D150").Value2 = 44806
$ws.Range("E150").Value2 = 9
$ws.Range("F150").Value2 = 100112013
$ws.Range("G150").Value2 = "Alcachofa"
$ws.Range("H150").Value2 = "Madrigal"
$ws.Range("I150").Value2 = "Primera"
$ws.Range("J150").Value2 = 65
$ws.Range("K150").Value2 = 13000
$ws.Range("L150").Value2 = 13000
$ws.Range("M150").Value2 = 13000
$ws.Range("N150").Value2 = "$/caja 40 unidades"
$ws.Range("O150").Value2 = "Provincia de Limarí"
$ws.Range("P150").Value2 = 325
$ws.Range("Q150").Value2 = 40
$ws.Range("R150").Value2 = "Hortaliza"
